$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1373.1
$ws.Range("I33").Value = 1171.0952
$ws.Range("K33").Value = 1171.0952
$ws.Range("M33").Value = -942.0952

$ws.Range("H92").Value = 11495868
$ws.Range("I92").Value = 15874436
$ws.Range("J92").Value = 2124.5
$ws.Range("K92").Value = 15874436
$ws.Range("L92").Value = 2124.5
$ws.Range("M92").Value = -15873188
$ws.Range("N92").Value = -4620.5

$ws.Range("H99").Value = 9017.691999999999
$ws.Range("I99").Value = 565.9
$ws.Range("J99").Value = 37190.332
$ws.Range("K99").Value = 1697.7
$ws.Range("L99").Value = 111570.996
$ws.Range("M99").Value = -199.6999999999998
$ws.Range("N99").Value = -114566.996

$ws.Range("H137").Value = 1195.4814
$ws.Range("I137").Value = 1080.8823
$ws.Range("J137").Value = 1390.3
$ws.Range("K137").Value = 3242.6469
$ws.Range("L137").Value = 4170.9
$ws.Range("M137").Value = -692.6468999999997
$ws.Range("N137").Value = -9270.9

$ws.Range("H138").Value = 2603.9805
$ws.Range("J138").Value = 2517.4324
$ws.Range("L138").Value = 7552.297200000001
$ws.Range("N138").Value = -17832.2972

$ws.Range("H141").Value = 5158.0435
$ws.Range("I141").Value = 2109.5
$ws.Range("K141").Value = 6328.5
$ws.Range("M141").Value = -1148.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2184.808
$ws.Range("I32").Value = 1708.4419
$ws.Range("J32").Value = 5336.154
$ws.Range("K32").Value = 1708.4419
$ws.Range("L32").Value = 5336.154
$ws.Range("M32").Value = -1421.4419
$ws.Range("N32").Value = -5910.154

$ws.Range("H61").Value = 7409619.5
$ws.Range("I61").Value = 19609208
$ws.Range("J61").Value = 2726.9285
$ws.Range("K61").Value = 19609208
$ws.Range("L61").Value = 2726.9285
$ws.Range("M61").Value = -19608996
$ws.Range("N61").Value = -3150.9285

$ws.Range("H74").Value = 866.46155
$ws.Range("I74").Value = 571.6667
$ws.Range("J74").Value = 1119.1428
$ws.Range("K74").Value = 571.6667
$ws.Range("L74").Value = 1119.1428
$ws.Range("M74").Value = 302.3333
$ws.Range("N74").Value = -2867.1428

$ws.Range("H77").Value = 866.46155
$ws.Range("I77").Value = 571.6667
$ws.Range("J77").Value = 1119.1428
$ws.Range("K77").Value = 2858.3335
$ws.Range("L77").Value = 5595.714
$ws.Range("M77").Value = 1509.6665
$ws.Range("N77").Value = -14331.714

$ws.Range("H132").Value = 5036.722
$ws.Range("I132").Value = 5600.4346
$ws.Range("K132").Value = 16801.3038
$ws.Range("M132").Value = -14271.3038

$ws.Range("H136").Value = 7409619.5
$ws.Range("I136").Value = 19609208
$ws.Range("J136").Value = 2726.9285
$ws.Range("K136").Value = 58827624
$ws.Range("L136").Value = 8180.7855
$ws.Range("M136").Value = -58825074
$ws.Range("N136").Value = -13280.7855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2515.7632
$ws.Range("I134").Value = 2197.9666
$ws.Range("K134").Value = 6593.899800000001
$ws.Range("M134").Value = -4058.899800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4586.55
$ws.Range("I31").Value = 1145.9615
$ws.Range("J31").Value = 10976.214
$ws.Range("K31").Value = 1145.9615
$ws.Range("L31").Value = 10976.214
$ws.Range("M31").Value = -850.9614999999999
$ws.Range("N31").Value = -11566.214

$ws.Range("H34").Value = 4586.55
$ws.Range("I34").Value = 1145.9615
$ws.Range("J34").Value = 10976.214
$ws.Range("K34").Value = 1145.9615
$ws.Range("L34").Value = 10976.214
$ws.Range("M34").Value = -943.9614999999999
$ws.Range("N34").Value = -11380.214

$ws.Range("H58").Value = 3161.077
$ws.Range("I58").Value = 3358
$ws.Range("J58").Value = 2504.6667
$ws.Range("K58").Value = 3358
$ws.Range("L58").Value = 2504.6667
$ws.Range("M58").Value = -3155
$ws.Range("N58").Value = -2910.6667

$ws.Range("H132").Value = 23813884
$ws.Range("I132").Value = 6006
$ws.Range("K132").Value = 18018
$ws.Range("M132").Value = -15488

$ws.Range("H134").Value = 1057.6875
$ws.Range("I134").Value = 724.8461
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 2174.5383
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = 360.4616999999998
$ws.Range("N134").Value = -12570

$ws.Range("H136").Value = 3161.077
$ws.Range("I136").Value = 3358
$ws.Range("J136").Value = 2504.6667
$ws.Range("K136").Value = 10074
$ws.Range("L136").Value = 7514.000100000001
$ws.Range("M136").Value = -7524
$ws.Range("N136").Value = -12614.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1449.3334
$ws.Range("I68").Value = 817.1429000000001
$ws.Range("J68").Value = 1660.0635
$ws.Range("K68").Value = 2451.4287
$ws.Range("L68").Value = 4980.1905
$ws.Range("M68").Value = -1640.4287
$ws.Range("N68").Value = -6602.1905

$ws.Range("H71").Value = 1449.3334
$ws.Range("I71").Value = 817.1429000000001
$ws.Range("J71").Value = 1660.0635
$ws.Range("K71").Value = 7354.2861
$ws.Range("L71").Value = 14940.5715
$ws.Range("M71").Value = -3298.2861
$ws.Range("N71").Value = -23052.5715

$ws.Range("H113").Value = 716.0172
$ws.Range("I113").Value = 438.8158
$ws.Range("K113").Value = 1316.4474
$ws.Range("M113").Value = 853.5526

$ws.Range("H141").Value = 11749.125
$ws.Range("I141").Value = 12826.667
$ws.Range("J141").Value = 8516.5
$ws.Range("K141").Value = 38480.001
$ws.Range("L141").Value = 25549.5
$ws.Range("M141").Value = -33300.001
$ws.Range("N141").Value = -35909.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H42").Value = 45374
$ws.Range("J42").Value = 45374
$ws.Range("L42").Value = 45374
$ws.Range("N42").Value = -46344

$ws.Range("H115").Value = 45374
$ws.Range("J115").Value = 45374
$ws.Range("L115").Value = 45374
$ws.Range("N115").Value = -47724

$ws.Range("H132").Value = 2440.8462
$ws.Range("I132").Value = 2022.0952
$ws.Range("K132").Value = 6066.2856
$ws.Range("M132").Value = -3536.2856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 70002
$ws.Range("J2").Value = 70002
$ws.Range("L2").Value = 70002
$ws.Range("N2").Value = -70226

$ws.Range("H22").Value = 9218.385
$ws.Range("I22").Value = 1100
$ws.Range("J22").Value = 14292.375
$ws.Range("K22").Value = 1100
$ws.Range("L22").Value = 14292.375
$ws.Range("M22").Value = -805
$ws.Range("N22").Value = -14882.375

$ws.Range("H27").Value = 9218.385
$ws.Range("I27").Value = 1100
$ws.Range("J27").Value = 14292.375
$ws.Range("K27").Value = 1100
$ws.Range("L27").Value = 14292.375
$ws.Range("M27").Value = -993
$ws.Range("N27").Value = -14506.375

$ws.Range("H61").Value = 2757.6191
$ws.Range("I61").Value = 2471.5386
$ws.Range("J61").Value = 3222.5
$ws.Range("K61").Value = 2471.5386
$ws.Range("L61").Value = 3222.5
$ws.Range("M61").Value = -2269.5386
$ws.Range("N61").Value = -3626.5

$ws.Range("H93").Value = 2064.6365
$ws.Range("J93").Value = 2501.5715
$ws.Range("L93").Value = 2501.5715
$ws.Range("N93").Value = -4997.5715

$ws.Range("H98").Value = 89285
$ws.Range("J98").Value = 89285
$ws.Range("L98").Value = 89285
$ws.Range("N98").Value = -95275

$ws.Range("H113").Value = 2757.6191
$ws.Range("I113").Value = 2471.5386
$ws.Range("J113").Value = 3222.5
$ws.Range("K113").Value = 2471.5386
$ws.Range("L113").Value = 3222.5
$ws.Range("M113").Value = -301.5385999999999
$ws.Range("N113").Value = -7562.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7677646.5
$ws.Range("I132").Value = 2198.65
$ws.Range("J132").Value = 16205922
$ws.Range("K132").Value = 6595.950000000001
$ws.Range("L132").Value = 48617766
$ws.Range("M132").Value = -4065.950000000001
$ws.Range("N132").Value = -48622826

$ws.Range("H136").Value = 3113.639
$ws.Range("I136").Value = 3024.8
$ws.Range("K136").Value = 9074.400000000001
$ws.Range("M136").Value = -6524.400000000001

$ws.Range("H137").Value = 54628.57
$ws.Range("J137").Value = 54628.57
$ws.Range("L137").Value = 54628.57
$ws.Range("N137").Value = -64828.57
